$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 265, shifting rows 265:373 down to 266:374.
$ws.Rows(265).Insert()

# Populate the newly inserted row 265 with the new record's data.
$ws.Cells.Item(265, 1).Value = 3
$ws.Cells.Item(265, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(265, 3).Value = "Coquimbo"
$ws.Cells.Item(265, 4).Value = Get-Date -Year 2022 -Month 10 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(265, 5).Value = 5
$ws.Cells.Item(265, 6).Value = 100112039
$ws.Cells.Item(265, 7).Value = "Ciboulette"
$ws.Cells.Item(265, 8).Value = "Sin especificar"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 83
$ws.Cells.Item(265, 11).Value = 1500
$ws.Cells.Item(265, 12).Value = 1800
$ws.Cells.Item(265, 13).Value = 1673
$ws.Cells.Item(265, 14).Value = "`$/docena de atados"
$ws.Cells.Item(265, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(265, 16).Value = 558
$ws.Cells.Item(265, 17).Value = 3
$ws.Cells.Item(265, 18).Value = "Hortaliza"
